# ZDD Bb4A Convertingtestcases.xlsx -- "Add files via upload" edit
#
# The underlying change is small:
#   * the tester-sign-off cell G18 on the "-the-name-of-your- module" sheet
#     is updated from the generic "????? Tester's Name ????? / ????? Date
#     ?????" placeholder to "Drashti desai / 11 August 2023";
#   * five new test-case rows (20-24) are filled in under the "Indexing
#     Strings" section, each with a scenario name, a +/- result, an actual
#     "quote" result, the full console-transcript text and a PASS verdict;
#   * the sheet's scroll position / selection is moved down to where the
#     user was last working (G17).
#
# Dropping the old placeholder string from the shared-string pool shifts
# every other shared-string index down by one -- that shift is an automatic
# side effect of the save, not something to reproduce by hand, so every
# other cell on the sheet is left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "-the-name-of-your- module"
$ws.Activate()

# --- New rows 20-24 ---------------------------------------------------
# Values are written in the same order the strings first appear in the
# finished workbook so the shared-string table comes out in the same
# sequence as the target file.

$ws.Range("D21").Value = "*** Start of Indexing Strings Demo ***                             Type an int numeric string (q - to quit):`n0`nConverted number is 0                                                                                                                          *** End of Indexing Strings Demo ***    "
$ws.Range("D20").Value = "Type an int numeric string (g - to quit):`nabcdef`nConverted number is 0`n"
$ws.Range("C20").Value = "`"abcdef`"`n"
$ws.Range("D23").Value = "Type an int numeric string (g - to quit):`n1.5578686`nConverted number is 1`n"
$ws.Range("C23").Value = "`"1.5578686`""
$ws.Range("D22").Value = "Type an int numeric string (q - to quit):`n0`nConverted number is 0`n"
$ws.Range("D24").Value = "Type an int numeric string (g - to quit):`n214743648`nConverted number is -214743648`n"
$ws.Range("C24").Value = "`"214743648`""

$ws.Range("A20").Value = "Handle exit input"
$ws.Range("B20").Value = "positive "
$ws.Range("F20").Value = "PASS"

$ws.Range("A21").Value = "Exiting with q"
$ws.Range("B21").Value = "positive "
$ws.Range("C21").Value = "`"q`""
$ws.Range("F21").Value = "PASS"

$ws.Range("A22").Value = "minimal edge case"
$ws.Range("B22").Value = "positive "
$ws.Range("C22").Value = "`"0`""
$ws.Range("F22").Value = "PASS"

$ws.Range("A23").Value = " typical case"
$ws.Range("B23").Value = "positive "
$ws.Range("F23").Value = "PASS"

$ws.Range("A24").Value = "maximal edge case"
$ws.Range("B24").Value = "positive "
$ws.Range("F24").Value = "PASS"

# Content now drives these rows' autofit height -- pin them to the same
# heights as the finished workbook.
$ws.Rows(20).RowHeight = 66
$ws.Rows(21).RowHeight = 105.6
$ws.Rows(22).RowHeight = 66
$ws.Rows(23).RowHeight = 66
$ws.Rows(24).RowHeight = 79.2

# --- Tester sign-off (G18) ---------------------------------------------
# Written last so it becomes the final entry appended to the shared-string
# table, matching the target file.
$ws.Range("G18").Value = "Drashti desai `n11 August 2023"

# --- Scroll position / selection ---------------------------------------
$ws.Range("G17").Select()
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 1
